$d = $word.ActiveDocument
$d.Content.Find.Execute("23×46=", $true, $false, $false, $false, $false, $true, 1, $false, "96×58=", 2) | Out-Null
$d.Content.Find.Execute("61×44=", $true, $false, $false, $false, $false, $true, 1, $false, "53×79=", 2) | Out-Null
$d.Content.Find.Execute("28×31=", $true, $false, $false, $false, $false, $true, 1, $false, "90×80=", 2) | Out-Null
$d.Content.Find.Execute("98×28=", $true, $false, $false, $false, $false, $true, 1, $false, "89×22=", 2) | Out-Null
$d.Content.Find.Execute("35×76=", $true, $false, $false, $false, $false, $true, 1, $false, "22×39=", 2) | Out-Null
$d.Content.Find.Execute("12×56=", $true, $false, $false, $false, $false, $true, 1, $false, "33×54=", 2) | Out-Null
$d.Content.Find.Execute("26×87=", $true, $false, $false, $false, $false, $true, 1, $false, "83×45=", 2) | Out-Null
$d.Content.Find.Execute("68×11=", $true, $false, $false, $false, $false, $true, 1, $false, "85×25=", 2) | Out-Null
$d.Content.Find.Execute("76×23=", $true, $false, $false, $false, $false, $true, 1, $false, "80×61=", 2) | Out-Null
$d.Content.Find.Execute("76×46=", $true, $false, $false, $false, $false, $true, 1, $false, "30×47=", 2) | Out-Null
$d.Content.Find.Execute("73×91=", $true, $false, $false, $false, $false, $true, 1, $false, "24×44=", 2) | Out-Null
$d.Content.Find.Execute("36×13=", $true, $false, $false, $false, $false, $true, 1, $false, "45×94=", 2) | Out-Null
$d.Content.Find.Execute("38×54=", $true, $false, $false, $false, $false, $true, 1, $false, "92×95=", 2) | Out-Null
$d.Content.Find.Execute("98×14=", $true, $false, $false, $false, $false, $true, 1, $false, "92×70=", 2) | Out-Null
$d.Content.Find.Execute("12×42=", $true, $false, $false, $false, $false, $true, 1, $false, "95×64=", 2) | Out-Null
$d.Content.Find.Execute("67×55=", $true, $false, $false, $false, $false, $true, 1, $false, "92×68=", 2) | Out-Null
$d.Content.Find.Execute("26×76=", $true, $false, $false, $false, $false, $true, 1, $false, "21×18=", 2) | Out-Null
$d.Content.Find.Execute("55×13=", $true, $false, $false, $false, $false, $true, 1, $false, "26×30=", 2) | Out-Null
$d.Content.Find.Execute("85×27=", $true, $false, $false, $false, $false, $true, 1, $false, "73×70=", 2) | Out-Null
$d.Content.Find.Execute("76×36=", $true, $false, $false, $false, $false, $true, 1, $false, "61×19=", 2) | Out-Null
$d.Content.Find.Execute("47×54=", $true, $false, $false, $false, $false, $true, 1, $false, "98×19=", 2) | Out-Null
$d.Content.Find.Execute("56×77=", $true, $false, $false, $false, $false, $true, 1, $false, "36×17=", 2) | Out-Null
$d.Content.Find.Execute("57×38=", $true, $false, $false, $false, $false, $true, 1, $false, "70×99=", 2) | Out-Null
$d.Content.Find.Execute("82×29=", $true, $false, $false, $false, $false, $true, 1, $false, "73×19=", 2) | Out-Null
$d.Content.Find.Execute("28×76=", $true, $false, $false, $false, $false, $true, 1, $false, "48×31=", 2) | Out-Null
